$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format-donor cell: untouched, default style (s=0), used to strip the
# temporary text-format style back off numeric-looking text cells after write.
$donor = $ws.Range("D4")

$ws.Range('D2').Value = '44.568.97'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '2.228.60'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.36'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '90.82'
$ws.Range('E6').Value = '  -4.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.560'
$ws.Range('E7').Value = '  -2.15%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  -4.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.85'
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('E11').Value = '  -3.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.98'
$ws.Range('E12').Value = '  -2.86%  '
$ws.Range('E13').Value = '  -0.50%  '
$ws.Range('D14').Value = '2.568.92'
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').Value = '2.332.66'
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.810'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.24'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D18').Value = '44.541.74'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('D19').Value = '0.0₃0909'
$ws.Range('E19').Value = '  -5.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.05'
$ws.Range('E20').Value = '  -5.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.38'
$ws.Range('E21').Value = '  -6.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.60'
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.23'
$ws.Range('E23').Value = '  -1.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.88'
$ws.Range('E24').Value = '  -7.67%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.92'
$ws.Range('E26').Value = '  -4.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.27'
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.53'
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.32'
$ws.Range('E29').Value = '  -8.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.61'
$ws.Range('E30').Value = '  -2.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.63'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '146.75'
$ws.Range('E32').Value = '  -4.28%  '
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0761'
$ws.Range('E34').Value = '  -4.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.05'
$ws.Range('E35').Value = '  -2.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.107'
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  -3.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.81'
$ws.Range('E38').Value = '  +3.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.52'
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.24'
$ws.Range('E40').Value = '  -7.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.66'
$ws.Range('E41').Value = '  -3.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0290'
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').Value = '1.782.71'
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.72'
$ws.Range('E45').Value = '  +7.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '79.73'
$ws.Range('E46').Value = '  -3.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.182'
$ws.Range('E47').Value = '  -5.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '96.06'
$ws.Range('E48').Value = '  -3.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.74'
$ws.Range('E49').Value = '  -3.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '66.86'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.67'
$ws.Range('E51').Value = '  -3.99%  '

# Re-copy the donors (default) formatting onto every cell we forced to
# text, so the text-format style does not linger as a new per-cell style.
$donor.Copy() | Out-Null
$ws.Range('D5,D6,D7,D9,D10,D11,D12,D16,D17,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D34,D35,D36,D37,D38,D39,D40,D41,D42,D45,D46,D47,D48,D49,D50,D51').PasteSpecial(-4122)
